# Fix PNAD 2009 "tentativa" data sheet:
# - column B header in row 2 should read "total" (was an erroneous
#   "unnamed: 1_level_1" placeholder left over from the pandas export)
# - the two section-header rows ("situação do domicílio" and
#   "grandes regiões e unidades da federação") have no data of their own;
#   they were pushed into the data block, leaving the "urbana"/"rural" and
#   region rows one row below where their numbers actually are. Removing
#   those two empty header rows re-aligns every label with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the column B sub-header in row 2.
$ws.Range("B2").Value = "total"

# Delete the "grandes regiões e unidades da federação" header row (row 8);
# deleting this first (bottom-most) keeps row 5's index valid for the
# next delete.
$ws.Rows.Item(8).Delete()

# Delete the "situação do domicílio" header row (row 5).
$ws.Rows.Item(5).Delete()
